$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (IROP)
$ws.Range("F2").Value = 142408862489.6619
$ws.Range("G2").Value = 123470149677.878
$ws.Range("H2").Value = 16701310436.44388
$ws.Range("I2").Value = 2237402375.34

# Row 9 (OP PIK)
$ws.Range("F9").Value = 142625904398.2969
$ws.Range("G9").Value = 67349336169.53
$ws.Range("I9").Value = 73948832770.09695

# Row 18 (OP VVV)
$ws.Range("F18").Value = 97147960741.49001
$ws.Range("G18").Value = 78302924392.54497
$ws.Range("H18").Value = 17029436532.725

# Row 19 (OP VVV)
$ws.Range("F19").Value = 360958057.4400001
$ws.Range("G19").Value = 255795847.375
$ws.Range("H19").Value = 105162210.065

# Row 23 (OP ZP)
$ws.Range("F23").Value = 35301809979.14999
$ws.Range("G23").Value = 29820859700.4
